$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the poly count value for the Shotgun feature asset row (D2)
$ws.Range("D2").Value = "300 Polygons"

# Update the selected/active cell in the sheet view
$ws.Range("D3").Select()
